# Rename the three header/footer logo pictures.
#
# The BTEC logo (in the header) and the two Pearson logos (one in each of
# the two footers) were exported with swapped image names - fix them up
# via InlineShape.Name, same as renaming a picture in the UI would do.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Footer 1 - Pearson logo: image2.png -> image1.png
$footer1 = $sec.Footers(1)
$footer1.Range.InlineShapes(1).Name = "image1.png"

# Footer 2 - Pearson logo: image2.png -> image1.png
$footer2 = $sec.Footers(2)
$footer2.Range.InlineShapes(1).Name = "image1.png"

# Header 2 - BTec logo: image1.jpg -> image2.jpg
$header2 = $sec.Headers(2)
$header2.Range.InlineShapes(1).Name = "image2.jpg"
